# Implemented getting number of lines for methods and classes.
# Adds two new worksheets:
#   - classNumberOfLines  : Class Name | Number of Lines
#   - methodNumberOfLines : Class Name | Method Signature | Number of Lines
# Both are appended after the existing "fieldInterfaceRelations" sheet.

$wb = $excel.ActiveWorkbook

# Helper: write a value as TEXT (shared string) even when it looks like a
# number ("6", "3", ...), instead of letting Excel auto-coerce it to a
# numeric cell. Briefly force a text number-format so the literal isn't
# parsed as a number, then clear the formatting back to the sheet default
# so the cell keeps using the workbook's default style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- classNumberOfLines sheet ---
# Worksheets.Add() inserts the new sheet at the front and makes it the
# ActiveSheet, and Move() always relocates whichever sheet is currently
# active - so populate + move this sheet to its final spot before adding
# the next one.
$wsClasses = $wb.Worksheets.Add()
$wsClasses.Name = "classNumberOfLines"
$wsClasses.Range("A1").Value = "Class Name"
$wsClasses.Range("B1").Value = "Number of Lines"
$wsClasses.Range("A2").Value = "org.andante.eureka.EurekaApplication"
Set-TextValue $wsClasses.Range("B2") "6"

$fieldInterfaceRelations = $wb.Worksheets.Item("fieldInterfaceRelations")
$wsClasses.Move($null, $fieldInterfaceRelations)

# --- methodNumberOfLines sheet ---
$wsMethods = $wb.Worksheets.Add()
$wsMethods.Name = "methodNumberOfLines"
$wsMethods.Range("A1").Value = "Class Name"
$wsMethods.Range("B1").Value = "Method Signature"
$wsMethods.Range("C1").Value = "Number of Lines"
$wsMethods.Range("A2").Value = "org.andante.eureka.EurekaApplication"
$wsMethods.Range("B2").Value = "main(java.lang.String[])"
Set-TextValue $wsMethods.Range("C2") "3"

# Re-resolve classNumberOfLines by name (the old handle's position is
# stale after the previous Move) and place methodNumberOfLines right
# after it, at the very end of the workbook.
$classesSheet = $wb.Worksheets.Item("classNumberOfLines")
$wsMethods.Move($null, $classesSheet)
